# Updated cryptos list (price + 1h volume change) per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '38.805.66'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').Value = '2.093.44'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '228.90'
$ws.Range('E5').Value = '  -0.21%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '61.64'
$ws.Range('E7').Value = '  +1.34%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.387'
$ws.Range('E9').Value = '  +1.79%  '
$ws.Range('E10').Value = '  +0.34%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.103'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '15.34'
$ws.Range('E12').Value = '  +4.78%  '
$ws.Range('D13').Value = '2.403.14'
$ws.Range('E13').Value = '  +0.03%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '22.12'
$ws.Range('E14').Value = '  -1.19%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.808'
$ws.Range('E15').Value = '  +4.26%  '
$ws.Range('E16').Value = '  -0.83%  '
$ws.Range('D17').Value = '2.070.40'
$ws.Range('E17').Value = '  -1.29%  '
$ws.Range('D18').Value = '38.766.56'
$ws.Range('E18').Value = '  +1.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '71.91'
$ws.Range('E19').Value = '  +2.40%  '
$ws.Range('E20').Value = '  +1.00%  '
$ws.Range('D21').Value = '0.0₃0840'
$ws.Range('E21').Value = '  +0.66%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '228.25'
$ws.Range('E22').Value = '  +1.74%  '
$ws.Range('E23').Value = '  -0.43%  '
$ws.Range('E24').Value = '  -2.50%  '
$ws.Range('E25').Value = '  +0.67%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '171.48'
$ws.Range('E26').Value = '  +0.75%  '
$ws.Range('E27').Value = '  +0.90%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.137'
$ws.Range('E28').Value = '  +4.52%  '
$ws.Range('E29').Value = '  +5.16%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '19.32'
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.51'
$ws.Range('E31').Value = '  +5.23%  '
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.52'
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.75'
$ws.Range('E34').Value = '  +1.40%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0618'
$ws.Range('E35').Value = '  +2.16%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.51'
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.39'
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('E38').Value = '  +1.57%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.20'
$ws.Range('E40').Value = '  +0.86%  '
$ws.Range('E41').Value = '  +4.09%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '101.30'
$ws.Range('E42').Value = '  +1.05%  '
$ws.Range('D43').Value = '1.534.26'
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.80'
$ws.Range('E44').Value = '  -1.28%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0911'
$ws.Range('E45').Value = '  +0.59%  '
$ws.Range('E46').Value = '  +2.06%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.66'
$ws.Range('E47').Value = '  +5.56%  '
$ws.Range('E48').Value = '  -1.03%  '
$ws.Range('E49').Value = '  +1.35%  '
$ws.Range('E50').Value = '  -1.01%  '
$ws.Range('D51').Value = '2.292.09'
$ws.Range('E51').Value = '  +0.11%  '

Write-Host "Updated cryptos list on $(Get-Date) with GitHub Actions"
